{"js": "// Map of old division expressions -> new division expressions.\nconst replacements = [\n  [\"821\u00f79=\", \"998\u00f72=\"],\n  [\"416\u00f78=\", \"802\u00f76=\"],\n  [\"760\u00f73=\", \"367\u00f72=\"],\n  [\"356\u00f78=\", \"740\u00f73=\"],\n  [\"665\u00f79=\", \"924\u00f73=\"],\n  [\"553\u00f79=\", \"115\u00f73=\"],\n  [\"901\u00f75=\", \"643\u00f76=\"],\n  [\"598\u00f77=\", \"179\u00f73=\"],\n  [\"976\u00f79=\", \"452\u00f78=\"],\n  [\"909\u00f76=\", \"527\u00f77=\"],\n  [\"686\u00f78=\", \"145\u00f77=\"],\n  [\"223\u00f77=\", \"858\u00f73=\"],\n  [\"637\u00f76=\", \"176\u00f74=\"],\n  [\"811\u00f73=\", \"767\u00f73=\"],\n  [\"172\u00f78=\", \"174\u00f72=\"],\n  [\"407\u00f73=\", \"636\u00f72=\"],\n  [\"412\u00f74=\", \"965\u00f72=\"],\n  [\"752\u00f76=\", \"176\u00f79=\"],\n  [\"349\u00f79=\", \"722\u00f79=\"],\n  [\"585\u00f77=\", \"664\u00f72=\"],\n  [\"218\u00f76=\", \"772\u00f73=\"],\n  [\"706\u00f79=\", \"241\u00f76=\"],\n  [\"542\u00f79=\", \"686\u00f72=\"],\n  [\"438\u00f78=\", \"606\u00f79=\"],\n  [\"850\u00f75=\", \"231\u00f76=\"],\n];\n\n// Use Find/Replace (whole-match, not wildcard) for each old->new pair.\n// MatchCase ensures exact textual replacement; these strings contain\n// no spaces so whole word matching is not needed/available for the\n// division sign anyway.\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Map of old division expressions -> new division expressions.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"821\u00f79=\", \"998\u00f72=\"),\n    @(\"416\u00f78=\", \"802\u00f76=\"),\n    @(\"760\u00f73=\", \"367\u00f72=\"),\n    @(\"356\u00f78=\", \"740\u00f73=\"),\n    @(\"665\u00f79=\", \"924\u00f73=\"),\n    @(\"553\u00f79=\", \"115\u00f73=\"),\n    @(\"901\u00f75=\", \"643\u00f76=\"),\n    @(\"598\u00f77=\", \"179\u00f73=\"),\n    @(\"976\u00f79=\", \"452\u00f78=\"),\n    @(\"909\u00f76=\", \"527\u00f77=\"),\n    @(\"686\u00f78=\", \"145\u00f77=\"),\n    @(\"223\u00f77=\", \"858\u00f73=\"),\n    @(\"637\u00f76=\", \"176\u00f74=\"),\n    @(\"811\u00f73=\", \"767\u00f73=\"),\n    @(\"172\u00f78=\", \"174\u00f72=\"),\n    @(\"407\u00f73=\", \"636\u00f72=\"),\n    @(\"412\u00f74=\", \"965\u00f72=\"),\n    @(\"752\u00f76=\", \"176\u00f79=\"),\n    @(\"349\u00f79=\", \"722\u00f79=\"),\n    @(\"585\u00f77=\", \"664\u00f72=\"),\n    @(\"218\u00f76=\", \"772\u00f73=\"),\n    @(\"706\u00f79=\", \"241\u00f76=\"),\n    @(\"542\u00f79=\", \"686\u00f72=\"),\n    @(\"438\u00f78=\", \"606\u00f79=\"),\n    @(\"850\u00f75=\", \"231\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n}\n"}
